# OPR293_DeliveryDocumentation_TestData.xlsx
# Commit: "Adding LTE001 and WHS001test data"
#
# Adds a new worksheet "OPR293_DLV_00006" (positioned between
# OPR293_DLV_00005 and OPR293_DLV_00010) populated with the same
# delivery-documentation row used on OPR293_DLV_00001, tweaks a few
# cell selections left behind by the editing session, and restores a
# portrait page setup on the first sheet.

$wb = $excel.ActiveWorkbook

# --- Insert the new sheet right before OPR293_DLV_00010 -------------------
# (Worksheets.Add() drops the new sheet immediately before the workbook's
# current ActiveSheet, which at load time is OPR293_DLV_00010.)
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "OPR293_DLV_00006"

# Re-resolve the other sheets *after* the insert so none of the handles are
# stale relative to the now-shifted worksheet collection.
$sheet1 = $wb.Worksheets.Item("OPR293_DLV_00001")
$sheet3 = $wb.Worksheets.Item("OPR293_DLV_00005")
$sheet10 = $wb.Worksheets.Item("OPR293_DLV_00010")

$headers = @("AgentCode", "ShipperCode ", "ConsigneeCode", "Origin", "Destination", `
    "ProductCode", "SCC", "Commodity", "ShipmentDescription", "ServiceCargoClass", `
    "Piece", "Weight", "ChargeType", "ModeOfPayment", "AWBSectionName", "cartType", `
    "Bdn_Locn", "Bdn_RcvdPieces", "Bdn_RcvdWeight")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$newSheet.Range("A2").Value = 11377
$newSheet.Range("B2").Value = 11377
$newSheet.Range("C2").Value = 11377
$newSheet.Range("D2").Value = "SFO"
$newSheet.Range("E2").Value = "LAX"
$newSheet.Range("F2").Value = "GOLDSTREAK"
$newSheet.Range("G2").Value = "None"
$newSheet.Range("H2").Value = "NONSCR"
$newSheet.Range("I2").Value = "None"
$newSheet.Range("J2").Value = "None"
$newSheet.Range("K2").Value = 2
$newSheet.Range("L2").Value = 59
$newSheet.Range("M2").Value = "CC"
$newSheet.Range("N2").Value = "None"
$newSheet.Range("O2").Value = "PlannedShipment"
$newSheet.Range("P2").Value = "CART"
$newSheet.Range("Q2").Value = "IDEFLOC"
$newSheet.Range("R2").Value = 2
$newSheet.Range("S2").Value = 59

# --- Restore a portrait page setup on OPR293_DLV_00001 ---------------------
$sheet1.PageSetup.Orientation = 1

# --- Leftover cell selections from the editing session ---------------------
[void]$sheet1.Range("C6").Select()
[void]$sheet3.Range("A1:XFD2").Select()
[void]$sheet10.Range("I9").Select()

# Leave the newly added sheet active/selected last, matching the saved
# workbook (it is the active tab, with its own leftover selection).
[void]$newSheet.Range("E6").Select()
